$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, centered, bordered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-65: new columns I (I0) and J (IF)
$data = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 8, 9),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 8, 8),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 9, 9),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 9, 9),
    @(38, 8, 8),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 8, 8),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 9, 9),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 9, 9),
    @(53, 8, 8),
    @(54, 8, 8),
    @(55, 8, 8),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 8, 8),
    @(59, 8, 8),
    @(60, 8, 8),
    @(61, 7, 7),
    @(62, 5, 5),
    @(63, 2, 2),
    @(64, 4, 4),
    @(65, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
